$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (grades) per the diff
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 5

$ws.Range("G12").Value = 5

$ws.Range("H29").Value = 5

# Update the frozen-pane scroll position and active selection to reflect
# where the author was working when the workbook was saved (pane stays
# split at the same row/column; only the visible top-left cell moves).
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 3
$ws.Range("A29").Select()
